# Refresh crypto price/volume data and re-order two swapped coin pairs,
# matching the upstream GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '66.679.42'
$ws.Range("E2").Value = '  +0.49%  '

# Row 3
$ws.Range("D3").Value = '3.244.88'
$ws.Range("E3").Value = '  +1.80%  '

# Row 4
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$ws.Range("E5").Value = '  +0.30%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.27'
$ws.Range("E6").Value = '  +0.85%  '

# Row 7
$ws.Range("E7").Value = '  -0.05%  '

# Row 8
$ws.Range("D8").Value = '3.246.19'
$ws.Range("E8").Value = '  +1.88%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.548'
$ws.Range("E9").Value = '  -0.09%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.162'
$ws.Range("E10").Value = '  +2.09%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.81'
$ws.Range("E11").Value = '  -0.96%  '

# Row 12
$ws.Range("E12").Value = '  -1.61%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000272'
$ws.Range("E13").Value = '  +2.79%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '39.12'
$ws.Range("E14").Value = '  +0.44%  '

# Row 15
$ws.Range("D15").Value = '3.779.35'
$ws.Range("E15").Value = '  +1.83%  '

# Row 16
$ws.Range("D16").Value = '66.706.06'
$ws.Range("E16").Value = '  +0.38%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.249.05'
$ws.Range("E17").Value = '  +1.91%  '

# Row 18
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.31'
$ws.Range("E18").Value = '  -1.13%  '

# Row 19
$ws.Range("E19").Value = '  +1.81%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '508.34'
$ws.Range("E20").Value = '  -1.07%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.36'
$ws.Range("E21").Value = '  -0.89%  '

# Row 22
$ws.Range("E22").Value = '  +1.57%  '

# Row 23
$ws.Range("E23").Value = '  -0.82%  '

# Row 24
$ws.Range("E24").Value = '  -1.51%  '

# Row 25
$ws.Range("B25").Value = 'Hedera'
$ws.Range("C25").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.178'
$ws.Range("E25").Value = '  +98.94%  '

# Row 26
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '86.17'
$ws.Range("E26").Value = '  +1.75%  '

# Row 27
$ws.Range("E27").Value = '  +0.12%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.02'
$ws.Range("E28").Value = '  +0.35%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.09'
$ws.Range("E29").Value = '  -1.14%  '

# Row 30
$ws.Range("E30").Value = '  -1.61%  '

# Row 31
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.91'
$ws.Range("E31").Value = '  -1.74%  '

# Row 32
$ws.Range("B32").Value = 'Stacks'
$ws.Range("C32").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.89'
$ws.Range("E32").Value = '  -5.96%  '

# Row 33
$ws.Range("E33").Value = '  +0.62%  '

# Row 34
$ws.Range("E34").Value = '  +0.01%  '

# Row 35
$ws.Range("E35").Value = '  -4.65%  '

# Row 36
$ws.Range("E36").Value = '  -2.69%  '

# Row 37
$ws.Range("D37").Value = '0.0₃0805'
$ws.Range("E37").Value = '  +19.34%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.30'
$ws.Range("E38").Value = '  +0.97%  '

# Row 39
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '493.94'
$ws.Range("E39").Value = '  -2.91%  '

# Row 40
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.28'
$ws.Range("E40").Value = '  +15.33%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0425'
$ws.Range("E41").Value = '  +0.55%  '

# Row 42
$ws.Range("E42").Value = '  +2.02%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.76'
$ws.Range("E43").Value = '  -1.22%  '

# Row 44
$ws.Range("E44").Value = '  -2.40%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.48'
$ws.Range("E45").Value = '  +1.64%  '

# Row 46
$ws.Range("D46").Value = '2.946.57'
$ws.Range("E46").Value = '  +3.28%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.33'
$ws.Range("E47").Value = '  -0.34%  '

# Row 48
$ws.Range("E48").Value = '  -0.10%  '

# Row 49
$ws.Range("E49").Value = '  +1.77%  '

# Row 50
$ws.Range("E50").Value = '  +0.02%  '

# Row 51
$ws.Range("E51").Value = '  -1.01%  '
